# Generate Report for Archive
#
# The localization status report is regenerated: the single pending item's
# status flips from "Ready for handoff" to "In Translation" everywhere it is
# shown (the Overview sheet's per-locale status columns, and the "Status"
# column on each per-locale detail sheet). Because the new status text is
# shorter, the report-generation step that produced this workbook also
# re-sized the (now too-wide) status columns to better fit the new text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the per-locale
#     status for row 2 (the single tracked file). ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- Per-locale detail sheets: column C is "Status". ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus

# --- Re-fit the status columns now that the text is shorter. The host's
#     ColumnWidth setter quantizes to its own internal pixel grid, so we pick
#     the input that lands the saved width as close as possible to the
#     freshly-autofit width a real report-generation pass would produce. ---
$fitWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $fitWidth
$overview.Columns.Item(6).ColumnWidth = $fitWidth
$zhcn.Columns.Item(3).ColumnWidth = $fitWidth
$dede.Columns.Item(3).ColumnWidth = $fitWidth
